$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Contacts page ("Places"): add a new row for the Nourishing Life Center
# ---------------------------------------------------------------------------
$places = $wb.Worksheets.Item("Places")

$places.Range("A14").Value = "NrshLifeCtr"
$places.Range("B14").Value = "Nourishing Life Center of Health"
$places.Range("C14").Value = "http://nourishinglife.com/yoga"
$places.Range("C14").Hyperlinks.Add($places.Range("C14"), "http://nourishinglife.com/yoga") | Out-Null
$places.Range("C14").Style = "Hyperlink"

$places.Range("A15").Select() | Out-Null

# ---------------------------------------------------------------------------
# 2. Add the new studio schedule sheet, "NrshLifeCtr", after "YmcaReut"
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "NrshLifeCtr"

# Column widths (approximate best-fit)
$ws.Columns.Item(1).ColumnWidth = 10.58
$ws.Columns.Item(4).ColumnWidth = 34.75
$ws.Columns.Item(5).ColumnWidth = 17.42

# Header row
$ws.Range("A1").Value = "DOW"
$ws.Range("B1").Value = "Start"
$ws.Range("C1").Value = "Stop"
$ws.Range("D1").Value = "What"
$ws.Range("E1").Value = "Who"

function Set-ClassRow($Row, $Day, $Start, $Stop, $What, $Who) {
    $ws.Range("A$Row").Value = $Day
    $ws.Range("B$Row").Value = $Start
    $ws.Range("B$Row").NumberFormat = "h:mm"
    $ws.Range("C$Row").Value = $Stop
    $ws.Range("C$Row").NumberFormat = "h:mm"
    $ws.Range("D$Row").Value = $What
    $ws.Range("E$Row").Value = $Who
}

Set-ClassRow 2  "Monday"    0.66666666666666663 0.71875             "Dao Flow All - Levels 1 & 2"          "Self Practice (sub)"
Set-ClassRow 4  "Tuesday"   0.39583333333333331 0.4375               "Gentle Flow"                          "Margaret Kirschner"
Set-ClassRow 6  "Wednesday" 0.72916666666666663 0.77083333333333337  "Dao Flow Yoga for Women Level 1"       "Dr. Robin Saraswati"
Set-ClassRow 8  "Friday"    0.39583333333333331 0.4375               "Dao Flow Yoga for Women Levels 1 & 2" "Self Practice (sub)"
Set-ClassRow 9  "Friday"    0.66666666666666663 0.71875              "Dao Flow All - Levels 1 & 2"          "Self Practice (sub)"
Set-ClassRow 11 "Saturday"  0.45833333333333331 0.52083333333333337  "Dao Flow Yoga for Women Level 2"       "Self Practice (sub)"

$ws.Range("C14").Select() | Out-Null

Write-Host "Edit complete"
